$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 8.2438
$ws.Range("A3").Value = -21.8031
$ws.Range("C3").Value = -11.2137
$ws.Range("E6").Value = 16.5157
$ws.Range("C12").Value = -12.04049999999999
$ws.Range("A14").Value = -21.6164
$ws.Range("E19").Value = 16.274
$ws.Range("A21").Value = -20.38609999999998
$ws.Range("A23").Value = -20.22339999999998
$ws.Range("C24").Value = -12.42449999999999
$ws.Range("E24").Value = 16.14859999999999
$ws.Range("A25").Value = -22.0495
$ws.Range("B25").Value = 5.571999999999996
$ws.Range("C25").Value = -12.89869999999999
$ws.Range("A26").Value = -21.14269999999997
$ws.Range("B27").Value = 6.141200000000005
$ws.Range("A29").Value = -20.60299999999999
$ws.Range("E30").Value = 15.80649999999999
$ws.Range("B31").Value = 5.400899999999999
$ws.Range("E31").Value = 16.62429999999999
$ws.Range("E33").Value = 16.79250000000001
$ws.Range("B39").Value = 9.573900000000007
$ws.Range("E42").Value = 16.5093
$ws.Range("B48").Value = 5.524100000000003
$ws.Range("C50").Value = -13.35539999999999
$ws.Range("B51").Value = 5.191300000000004
$ws.Range("B52").Value = 4.700200000000002
$ws.Range("A53").Value = -21.55210000000001
$ws.Range("C53").Value = -10.3708
$ws.Range("B55").Value = 6.213299999999998
$ws.Range("E55").Value = 16.4255
$ws.Range("B56").Value = 4.6976
$ws.Range("A57").Value = -22.0542
$ws.Range("B57").Value = 5.1033
$ws.Range("C57").Value = -12.70319999999999
$ws.Range("E58").Value = 16.14300000000002
$ws.Range("A59").Value = -22.24970000000001
$ws.Range("C61").Value = -13.19009999999999
$ws.Range("C63").Value = -12.1299
$ws.Range("E65").Value = 16.82850000000001
$ws.Range("A69").Value = -21.62229999999998
$ws.Range("C70").Value = -12.5612
$ws.Range("E70").Value = 16.86710000000001
$ws.Range("B73").Value = 8.394699999999998
$ws.Range("E75").Value = 16.35130000000001
$ws.Range("A79").Value = -20.5483
$ws.Range("A83").Value = -22.0459
$ws.Range("E83").Value = 16.2133
$ws.Range("C86").Value = -13.2638
$ws.Range("E86").Value = 16.70320000000001
$ws.Range("B89").Value = 5.156499999999995
$ws.Range("B90").Value = 5.509600000000002
$ws.Range("A91").Value = -21.2933
$ws.Range("B92").Value = 5.233199999999993
$ws.Range("A93").Value = -20.75439999999999
$ws.Range("E96").Value = 16.14279999999999
$ws.Range("E97").Value = 16.67890000000001
$ws.Range("C98").Value = -11.41839999999999
$ws.Range("C100").Value = -12.44259999999999
$ws.Range("C102").Value = -13.4799
